$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "{'correlação_pearson': np.float64(-0.02), 'p_valor': np.float64(0.24)}"
$ws.Range("B3").Value = "{'qui2': np.float64(0.55), 'p_valor': np.float64(0.46), 'graus_liberdade': 1, 'tabela': {1.0: {1.0: 1500, 2.0: 1437}, 2.0: {1.0: 1534, 2.0: 1529}}}"
$ws.Range("B4").Value = "{'qui2': np.float64(3.01), 'p_valor': np.float64(0.08), 'graus_liberdade': 1, 'tabela': {1.0: {1.0: 1502, 2.0: 1504}, 2.0: {1.0: 1428, 2.0: 1566}}}"
$ws.Range("B5").Value = "{1.0: 5.48, 2.0: 5.59, 3.0: 5.5}"
$ws.Range("B6").Value = "{0: {'Q5': 1.0, 'Q6': 1.51, 'Q7': 4.0, 'Q8': 1.99, 'Q9': 14.06, 'Q10': 1.51, 'Q11': 3.92, 'Q12': 2.01, 'Q13': 5.54, 'Q14': 1.52, 'Q15': 1.5, 'Q16': 1.5, 'Q17': 1.52, 'Q18': 1.48, 'Q19': 1.49, 'Q20': 1.51, 'Q21': 1.49, 'Q22': 1.51, 'Q23': 1.5, 'Q24': 1.52, 'Q25': 1.52, 'Q26': 2.01}, 1: {'Q5': 2.0, 'Q6': 1.45, 'Q7': 4.31, 'Q8': 2.01, 'Q9': 13.27, 'Q10': 1.49, 'Q11': 3.93, 'Q12': 1.31, 'Q13': 5.29, 'Q14': 1.54, 'Q15': 1.51, 'Q16': 1.44, 'Q17': 1.53, 'Q18': 1.49, 'Q19': 1.41, 'Q20': 1.5, 'Q21': 1.51, 'Q22': 1.55, 'Q23': 1.47, 'Q24': 1.46, 'Q25': 1.47, 'Q26': 2.02}, 2: {'Q5': 2.0, 'Q6': 1.52, 'Q7': 3.73, 'Q8': 1.99, 'Q9': 14.99, 'Q10': 1.49, 'Q11': 4.07, 'Q12': 2.64, 'Q13': 5.73, 'Q14': 1.47, 'Q15': 1.52, 'Q16': 1.55, 'Q17': 1.49, 'Q18': 1.49, 'Q19': 1.62, 'Q20': 1.5, 'Q21': 1.51, 'Q22': 1.46, 'Q23': 1.56, 'Q24': 1.56, 'Q25': 1.54, 'Q26': 1.99}}"
